$d = $word.ActiveDocument

# NOTE: the upstream commit's styles.xml hunk only adds four built-in
# <w:lsdException> rows ("Normal Table", "Table Subtle 1", "Table Web 2",
# "Table Web 3") to the latent-styles table so it matches the w:count="377"
# already declared on <w:latentStyles>. That table is Word's own style-
# gallery bookkeeping; it isn't backed by any Styles/Style property on the
# object model (there is no Styles.Add / latent-style surface to script),
# so it is left alone here - every other line of the diff is a genuine
# content edit and is reproduced below.

# ---------------------------------------------------------------------------
# 1) "This program requires Java to run. To execute the program simply run
#    the executable file." -> "To execute the program simply run the
#    executable file."
#
#    The real Word edit also leaves the "_GoBack" last-edit bookmark right
#    after the now-shortened sentence (it used to sit after the "Main"
#    heading run instead). We first park a one-character placeholder at the
#    very end of that paragraph, drop a bookmark immediately in front of it
#    (placing a bookmark collapsed exactly at "end of paragraph" tends to
#    slide around formatting-run boundaries, so nudging off that edge with a
#    throwaway character keeps the insertion point stable), then remove the
#    placeholder again. Adding a bookmark literally named "_GoBack" replaces
#    any existing "_GoBack" bookmark elsewhere in the story, so this single
#    call both plants the new one and removes the old one in front of "Main".
# ---------------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("This program requires Java to run. To execute the program simply run the executable file.")
if ($found) {
    $para = $target.Paragraphs(1)
    $paraEnd = $para.Range.End

    $anchor = $d.Range($paraEnd - 1, $paraEnd - 1)
    $anchor.InsertAfter("~")

    $bookmarkSpot = $d.Range($paraEnd - 1, $paraEnd - 1)
    $d.Bookmarks.Add("_GoBack", $bookmarkSpot)

    $placeholder = $d.Range($paraEnd - 1, $paraEnd)
    $placeholder.Delete()
}

$d.Content.Find.Execute("This program requires Java to run. To execute the program simply run the executable file.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "To execute the program simply run the executable file.", 2)

# ---------------------------------------------------------------------------
# 2) Collapse the three split runs in the table cell ("Determines the " +
#    "Column " + "that a piece of the diagonal line will be places on")
#    into the single run "Determines the Column that a piece of the
#    diagonal line will be places on".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Determines the Column that a piece of the diagonal line will be places on", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Determines the Column that a piece of the diagonal line will be places on", 2)
